$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Natalie's - Honey Tangerine (Quantity 2 -> 1, Total Cost 28.00 -> 14.00)
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "1"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "14.00"

# Row 7: Natalie's - Lemonade (Quantity 2 -> 1, Total Cost 18.50 -> 9.25)
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "1"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "9.25"

# Row 8: Natalie's - Orange Juice (Quantity 5 -> 3, Total Cost 118.75 -> 71.25)
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "3"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "71.25"

# Row 9: Natalie's - Orange Mango (Quantity 2 -> 1, Total Cost 26.00 -> 13.00)
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "1"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "13.00"

# Row 10: Natalie's - Orange Pineapple (Quantity 2 -> 1, Total Cost 26.00 -> 13.00)
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "1"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "13.00"
